$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) values must remain text, even when they look like plain
# numbers (e.g. "0.999"), to match the original inline-string cell type.
# Forcing NumberFormat="@" then resetting Style="Normal" keeps the cell on
# the default style (no stray "s" attribute) while the stored value stays text.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.216.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.316.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.652"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +6.03%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.314.45"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("E10").Value = "  -1.00%  "
$ws.Range("E11").Value = "  +2.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.401"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.893.19"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("E14").Value = "  -2.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "66.229.06"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.371.25"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000164"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "424.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.66%  "
$ws.Range("E20").Value = "  -2.83%  "
$ws.Range("E21").Value = "  -3.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.69"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.473.67"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.513"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.17%  "
$ws.Range("E28").Value = "  +6.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000115"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.28%  "
$ws.Range("E31").Value = "  +0.22%  "
$ws.Range("E32").Value = "  -1.60%  "
$ws.Range("E33").Value = "  -1.66%  "
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("E35").Value = "  -1.24%  "
$ws.Range("E36").Value = "  -1.95%  "
$ws.Range("E37").Value = "  -2.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "159.91"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.49%  "
$ws.Range("E39").Value = "  -2.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.861.19"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.88%  "
$ws.Range("E41").Value = "  +0.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.32"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.31%  "
$ws.Range("E43").Value = "  -2.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.758"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "39.69"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0660"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.41%  "
$ws.Range("E48").Value = "  -1.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.07"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "310.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.20%  "
$ws.Range("E51").Value = "  +0.44%  "
